$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "bkhbKHR"
$ws.Range("B6").Value = "BCJKSBDKVB;"
$ws.Range("C6").Select()
